# Fruta / hortaliza, semanal
# Insert a new weekly price record for Espárragos (Macroferia Regional de Talca)
# as row 89, pushing the existing rows 89-96 down to 90-97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 89 (shifts rows 89:96 -> 90:97)
$ws.Rows.Item(89).Insert()

# Populate the new row 89 with the latest survey data
$ws.Range("A89").Value = 5
$ws.Range("B89").Value = "Macroferia Regional de Talca"
$ws.Range("C89").Value = "Maule"
$ws.Range("D89").Value = 45173
$ws.Range("E89").Value = 7
$ws.Range("F89").Value = 300000000
$ws.Range("G89").Value = "Espárragos"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 2500
$ws.Range("L89").Value = 2500
$ws.Range("M89").Value = 2500
$ws.Range("N89").Value = "$/kilo"
$ws.Range("O89").Value = "Provincia de Linares"
$ws.Range("P89").Value = 2500
$ws.Range("Q89").Value = 1
$ws.Range("R89").Value = "Hortaliza"
